$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.121.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3900"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07891"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9898"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.064"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.763"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06995"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009988"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.118.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.318"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.113.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.109"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.955"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  -5.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09320"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8991"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.251"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.76%  "

$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05787"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.177"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("E41").Value = "  -0.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1790"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.709"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.170"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07013"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.854"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.557"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.045"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
